# Atualização de bases das ligas, do dia: 30-05-2024 às 23:16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header renames: ht_goals_h -> HTHG, ht_goals_a -> HTAG
$ws.Range("I1").Value = "HTHG"
$ws.Range("J1").Value = "HTAG"

# 2) Rows 11 and 12 swap their data (everything except column A / id)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$row11vals = @{}
$row12vals = @{}
foreach ($c in $cols) {
    $row11vals[$c] = $ws.Range("$c" + "11").Value2
    $row12vals[$c] = $ws.Range("$c" + "12").Value2
}

foreach ($c in $cols) {
    $ws.Range("$c" + "11").Value = $row12vals[$c]
    $ws.Range("$c" + "12").Value = $row11vals[$c]
}

# 3) Row 118 odds updates
$ws.Range("O118").Value = 1.666
$ws.Range("P118").Value = 3.5
$ws.Range("S118").Value = 1.85
$ws.Range("T118").Value = 1.95
$ws.Range("U118").Value = 2.5
$ws.Range("V118").Value = 1.825
$ws.Range("W118").Value = 1.975
